# Insert a new weekly record at row 327, shifting all existing rows
# (327-392) down by one (to 328-393), and fill the new row 327 with the
# new data point. This mirrors the target diff: dimension grows from
# A1:R392 to A1:R393, and a new record is inserted into the middle of
# the data block (all subsequent rows shift down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 327; this shifts rows 327:392 down
# to 328:393 and extends the used range / dimension automatically.
$ws.Rows.Item(327).Insert()

# Populate the newly inserted row 327 with the new data values.
$ws.Cells.Item(327, 1).Value2 = 9
$ws.Cells.Item(327, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(327, 3).Value2 = "Metropolitana"
$ws.Cells.Item(327, 4).Value2 = 45173
$ws.Cells.Item(327, 5).Value2 = 13
$ws.Cells.Item(327, 6).Value2 = 100112001
$ws.Cells.Item(327, 7).Value2 = "Berenjena"
$ws.Cells.Item(327, 8).Value2 = "Sin especificar"
$ws.Cells.Item(327, 9).Value2 = "Primera"
$ws.Cells.Item(327, 10).Value2 = 124
$ws.Cells.Item(327, 11).Value2 = 7000
$ws.Cells.Item(327, 12).Value2 = 8000
$ws.Cells.Item(327, 13).Value2 = 7500
$ws.Cells.Item(327, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(327, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(327, 16).Value2 = 150
$ws.Cells.Item(327, 17).Value2 = 50
$ws.Cells.Item(327, 18).Value2 = "Hortaliza"
